$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.909.19"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.588.14"
$ws.Range("E3").Value = "  -1.12%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "210.63"
$r.ClearFormats()
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("E8").Value = "  +0.98%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.0615"
$r.ClearFormats()
$ws.Range("E9").Value = "  -0.13%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "18.24"
$r.ClearFormats()
$ws.Range("E10").Value = "  +1.39%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0790"
$r.ClearFormats()
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "1.807.99"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").Value = "1.578.07"
$ws.Range("E13").Value = "  -1.79%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "4.04"
$r.ClearFormats()
$ws.Range("E14").Value = "  -0.84%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.506"
$r.ClearFormats()
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").Value = "25.896.68"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "0.0₃0726"
$ws.Range("E17").Value = "  -0.35%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "60.22"
$r.ClearFormats()
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("E19").Value = "  -0.44%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "194.93"
$r.ClearFormats()
$ws.Range("E20").Value = "  +2.56%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "4.21"
$r.ClearFormats()
$ws.Range("E21").Value = "  +0.34%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "9.38"
$r.ClearFormats()
$ws.Range("E22").Value = "  +0.32%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "5.96"
$r.ClearFormats()
$ws.Range("E23").Value = "  +0.14%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "0.131"
$r.ClearFormats()
$ws.Range("E24").Value = "  +1.01%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "141.15"
$r.ClearFormats()
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("E26").Value = "  -0.51%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "1.71"
$r.ClearFormats()
$ws.Range("E27").Value = "  -0.02%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "15.11"
$r.ClearFormats()
$ws.Range("E28").Value = "  +0.84%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "6.47"
$r.ClearFormats()
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("E30").Value = "  -3.92%  "
$ws.Range("E31").Value = "  +0.45%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "3.14"
$r.ClearFormats()
$ws.Range("E32").Value = "  +1.44%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "3.03"
$r.ClearFormats()
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("D36").Value = "1.097.96"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("E37").Value = "  -0.43%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "2.34"
$r.ClearFormats()
$ws.Range("E38").Value = "  -1.84%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.0151"
$r.ClearFormats()
$ws.Range("E39").Value = "  +0.16%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.503"
$r.ClearFormats()
$ws.Range("E40").Value = "  +0.51%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.784"
$r.ClearFormats()
$ws.Range("E41").Value = "  -3.32%  "
$ws.Range("E42").Value = "  +6.97%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "93.02"
$r.ClearFormats()
$ws.Range("E43").Value = "  -3.21%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "5.11"
$r.ClearFormats()
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("D45").Value = "1.720.02"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("E47").Value = "  +3.41%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "53.31"
$r.ClearFormats()
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("E51").Value = "  -0.56%  "
